$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bayern Munich - Club Brugge KV
$ws.Range("C2").Value = 74
$ws.Range("F2").Value = 1.22

# Row 3: Galatasaray -> Real Madrid vs Juventus FC
$ws.Range("A3").Value = "Real Madrid  - Juventus FC: 20:00"
$ws.Range("B3").Value = "Real Madrid"
$ws.Range("D3").Value = 94

# Row 4: Chelsea FC - Ajax Amsterdam
$ws.Range("C4").Value = 70
$ws.Range("D4").Value = 94
$ws.Range("F4").Value = 1.3
